# Applies the "updated cryptos list" data refresh described by the commit
# message / XML diff: per-row price + 1h-volume% updates, plus a 3-row
# reorder among Monero / PancakeSwap / Aptos (rows 29-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-ambiguous values ---
$ws.Range("D2").Value = '58.007.81'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '2.468.04'
$ws.Range("E3").Value = '  -2.33%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  -3.62%  '
$ws.Range("E6").Value = '  -4.68%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -2.38%  '
$ws.Range("E9").Value = '  -2.31%  '
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("E12").Value = '  -1.67%  '
$ws.Range("D13").Value = '2.906.35'
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("D14").Value = '57.934.20'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("E15").Value = '  -3.60%  '
$ws.Range("E16").Value = '  -2.47%  '
$ws.Range("D17").Value = '2.464.45'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("E22").Value = '  -3.66%  '
$ws.Range("E23").Value = '  -2.67%  '
$ws.Range("E24").Value = '  -3.21%  '
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  -3.54%  '
$ws.Range("E27").Value = '  -3.11%  '
$ws.Range("E28").Value = '  -3.18%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E30").Value = '  -4.71%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("E31").Value = '  -6.15%  '
$ws.Range("E32").Value = '  -2.27%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("E36").Value = '  -10.87%  '
$ws.Range("E37").Value = '  -3.87%  '
$ws.Range("E38").Value = '  -4.89%  '
$ws.Range("E39").Value = '  -3.17%  '
$ws.Range("E40").Value = '  -4.70%  '
$ws.Range("E41").Value = '  -4.64%  '
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("E44").Value = '  -4.95%  '
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("E46").Value = '  -4.18%  '
$ws.Range("E47").Value = '  -3.28%  '
$ws.Range("E48").Value = '  -2.09%  '
$ws.Range("D49").Value = '1.732.25'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("E50").Value = '  -1.42%  '

# --- Numeric-looking price values that must stay text (force text, then restore default style) ---
$numericTextCells = @("D5", "D6", "D7", "D15", "D18", "D20", "D21", "D22", "D23", "D25", "D29", "D30", "D31", "D35", "D37", "D39", "D41", "D42", "D43", "D44", "D48")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '518.03'
$ws.Range("D6").Value = '130.79'
$ws.Range("D7").Value = '1.00'
$ws.Range("D15").Value = '22.27'
$ws.Range("D18").Value = '10.78'
$ws.Range("D20").Value = '318.63'
$ws.Range("D21").Value = '1.00'
$ws.Range("D22").Value = '5.73'
$ws.Range("D23").Value = '64.07'
$ws.Range("D25").Value = '1.00'
$ws.Range("D29").Value = '165.89'
$ws.Range("D30").Value = '1.69'
$ws.Range("D31").Value = '6.30'
$ws.Range("D35").Value = '18.02'
$ws.Range("D37").Value = '3.97'
$ws.Range("D39").Value = '0.789'
$ws.Range("D41").Value = '272.11'
$ws.Range("D42").Value = '5.00'
$ws.Range("D43").Value = '0.592'
$ws.Range("D44").Value = '126.27'
$ws.Range("D48").Value = '17.03'

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
